# Update the marksheet's "Corr/total marks" figures on the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct-answer count goes from 3 -> 5
$ws.Range("B11").Value = 5

# Total row: total marks go from 63 -> 105
$ws.Range("B12").Value = 105

# Total row: the "obtained/max" summary text goes from 62/84 -> 105/140
$ws.Range("E12").Value = "105/140"
